$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: tiny floating point precision correction on the timestamp value
$ws.Range("A10").Value = 45877.37517245371

# Row 11: new data row appended
$ws.Range("A11").Value = 45877.41691968337
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 15.06
$ws.Range("E11").Value = 89.66
$ws.Range("F11").Value = 186.33
$ws.Range("G11").Value = 10.46
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "10:00:21"
